$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '72.753.50'
$ws.Range('E2').Value = '  -0.69%  '
$ws.Range('D3').Value = '3.950.66'
$ws.Range('E3').Value = '  -2.50%  '
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '606.86'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  +1.56%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '171.03'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  +11.14%  '
$ws.Range('E7').Value = '  -1.01%  '
$ws.Range('E8').Value = '  -0.02%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.787'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  +3.48%  '
$ws.Range('E10').Value = '  +5.88%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '56.02'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  +4.21%  '
$ws.Range('E12').Value = '  +0.67%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '11.54'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '  +4.33%  '
$ws.Range('D14').Value = '4.584.28'
$ws.Range('E14').Value = '  -2.55%  '
$ws.Range('D15').Value = '3.958.12'
$ws.Range('E15').Value = '  -2.36%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '21.38'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '  +2.81%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '14.05'
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').Value = '  -1.57%  '
$ws.Range('E18').Value = '  -1.26%  '
$ws.Range('B19').Value = 'WrappedBTC'
$ws.Range('C19').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D19').Value = '72.710.56'
$ws.Range('E19').Value = '  -0.72%  '
$ws.Range('B20').Value = 'TRON'
$ws.Range('C20').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '0.131'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  -0.82%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '445.05'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  +0.60%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '4.86'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  +0.49%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '95.68'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  -2.05%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '3.35'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  -5.08%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '14.23'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  -1.91%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '4.25'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  -2.51%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '11.33'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  -0.79%  '
$ws.Range('E28').Value = '  -1.14%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '10.40'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  -4.44%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '35.87'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  -3.12%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '7.95'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  -0.07%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '13.91'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  +1.83%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '49.47'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  +0.56%  '
$ws.Range('E34').Value = '  -4.37%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '0.0000100'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  +14.16%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '69.02'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  -3.29%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '631.92'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  -8.23%  '
$ws.Range('E38').Value = '  -3.99%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '3.47'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  +2.88%  '
$ws.Range('B40').Value = 'Dai'
$ws.Range('C40').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.999'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  +0.02%  '
$ws.Range('B41').Value = 'Kaspa'
$ws.Range('C41').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.146'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  -1.39%  '
$ws.Range('E42').Value = '  +0.09%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.0479'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  -3.60%  '
$ws.Range('E44').Value = '  -6.38%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '3.17'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  +42.19%  '
$ws.Range('E46').Value = '  -2.23%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '2.64'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  -2.23%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '3.38'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  -0.58%  '
$ws.Range('B49').Value = 'FLOKI'
$ws.Range('C49').Value = 'https://coinranking.com/coin/fmHk13Rqw+floki-floki'
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '0.000286'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  +5.78%  '
$ws.Range('B50').Value = 'WEMIXToken'
$ws.Range('C50').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '2.83'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  -16.26%  '
$ws.Range('D51').Value = '2.828.51'
$ws.Range('E51').Value = '  +1.16%  '
